$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @("09/13/2025", "DeepSeekV3", 4, "Soil", 12, 61, 35, 0.1643835616438356, 0.2553191489361702, 0.2, 47, 18, 55, 29, 0.2465753424657534, 0.3829787234042553, 0.3, 47),
        @("09/13/2025", "DeepSeekV3", 4, "city", 9, 9, 8, 0.5, 0.5294117647058824, 0.5142857142857143, 17, 9, 9, 8, 0.5, 0.5294117647058824, 0.5142857142857143, 17),
        @("09/13/2025", "DeepSeekV3", 4, "country", 9, 7, 13, 0.5625, 0.4090909090909091, 0.4736842105263158, 22, 9, 7, 13, 0.5625, 0.4090909090909091, 0.4736842105263158, 22),
        @("09/13/2025", "DeepSeekV3", 4, "cropSpecies", 30, 20, 51, 0.6, 0.3703703703703703, 0.4580152671755725, 81, 37, 13, 44, 0.74, 0.4567901234567901, 0.564885496183206, 81),
        @("09/13/2025", "DeepSeekV3", 4, "duration", 3, 1, 19, 0.75, 0.1363636363636364, 0.2307692307692308, 22, 4, 0, 18, 1, 0.1818181818181818, 0.3076923076923077, 22),
        @("09/13/2025", "DeepSeekV3", 4, "endTime", 20, 14, 11, 0.5882352941176471, 0.6451612903225806, 0.6153846153846154, 31, 21, 13, 10, 0.6176470588235294, 0.6774193548387096, 0.6461538461538462, 31),
        @("09/13/2025", "DeepSeekV3", 4, "region", 3, 3, 9, 0.5, 0.25, 0.3333333333333333, 12, 3, 3, 9, 0.5, 0.25, 0.3333333333333333, 12),
        @("09/13/2025", "DeepSeekV3", 4, "startTime", 31, 4, 26, 0.8857142857142857, 0.543859649122807, 0.673913043478261, 57, 32, 3, 25, 0.9142857142857143, 0.5614035087719298, 0.6956521739130436, 57)
)

$startRow = 10
$endRow = $startRow + $data.Length - 1

# Force column A (dates stored as literal text, e.g. "09/13/2025") to Text
# format so Excel does not auto-convert the strings into date serial
# numbers, then restore the default "Normal" style so the cells end up
# looking just like the existing rows (no explicit style attribute).
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = $rowValues[$col - 1]
    }
}

$ws.Range("A$startRow`:A$endRow").Style = "Normal"
